## Add 2022-Q4 data (feat: add 2022-Q4 data)
##
## 1. Update the "总计" (summary) sheet: insert a new top data row for
##    2022-Q4 and shift the existing quarters down by one row.
## 2. Insert a brand-new "2022-Q4" worksheet (right after "总计", i.e.
##    before the existing "2022-Q3" sheet) holding the three fund rows
##    for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet - rewrite the data rows (row 1 header is kept)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @("2022-Q4", 3, 0.03),
    @("2022-Q3", 1, 0.1),
    @("2022-Q2", 4, 0.24),
    @("2022-Q1", 5, 0.53),
    @("2021-Q4", 1, 3.36),
    @("2021-Q1", 1, 0),
    @("2020-Q4", 2, 1.08)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q4" sheet, inserted before the existing "2022-Q3"
#    tab so the tab order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, ...
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3Sheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$funds = @(
    @("519615", "银河君尚灵活配置混合I", "1.83", "38.98", "0.88", "0.0161", 8),
    @("519613", "银河君尚灵活配置混合A", "1.17", "38.98", "0.88", "0.0103", 8),
    @("519614", "银河君尚灵活配置混合C", "0.16", "38.98", "0.88", "0.0014", 8)
)

for ($i = 0; $i -lt $funds.Count; $i++) {
    $r = $i + 2
    $f = $funds[$i]

    $q4.Cells.Item($r, 1).Value = $i

    # Fund code must stay textual (it has no intrinsic numeric meaning and
    # some codes have leading zeros), so force text formatting before/after
    # the assignment instead of leaving it to auto-detection.
    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $f[0]
    $q4.Cells.Item($r, 2).Style = "Normal"

    $q4.Cells.Item($r, 3).Value = $f[1]
    $q4.Cells.Item($r, 4).Value = $f[2]
    $q4.Cells.Item($r, 5).Value = $f[3]
    $q4.Cells.Item($r, 6).Value = $f[4]
    $q4.Cells.Item($r, 7).Value = $f[5]
    $q4.Cells.Item($r, 8).Value = $f[6]
}
